# Reorder GO-term enrichment rows across sheets to reflect the refreshed
# enrichment run that now also compares the highest-CAI / lowest-CAI gene
# subsets (AT-rich / GC-rich sheets reorder a handful of terms; the
# Least-adapted / Most-adapted sheets are rewritten unchanged so the whole
# shared-string table is rebuilt in the same scan order Excel itself uses).

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("AT-rich genes")
$ws1.Range("A2:C34").ClearContents()
$ws2 = $wb.Worksheets.Item("GC-rich genes")
$ws2.Range("A2:C29").ClearContents()
$ws3 = $wb.Worksheets.Item("Least-adapted genes")
$ws3.Range("A2:C5").ClearContents()
$ws4 = $wb.Worksheets.Item("Most-adapted genes")
$ws4.Range("A2:C7").ClearContents()

# --- AT-rich genes ---
$ws1.Range("A2").Value = "GO:0030204"
$ws1.Range("B2").Value = "GO:BP"
$ws1.Range("C2").Value = "chondroitin sulfate metabolic process"
$ws1.Range("A3").Value = "GO:0030206"
$ws1.Range("B3").Value = "GO:BP"
$ws1.Range("C3").Value = "chondroitin sulfate biosynthetic process"
$ws1.Range("A4").Value = "GO:0050650"
$ws1.Range("B4").Value = "GO:BP"
$ws1.Range("C4").Value = "chondroitin sulfate proteoglycan biosynthetic process"
$ws1.Range("A5").Value = "GO:1903510"
$ws1.Range("B5").Value = "GO:BP"
$ws1.Range("C5").Value = "mucopolysaccharide metabolic process"
$ws1.Range("A6").Value = "GO:0050654"
$ws1.Range("B6").Value = "GO:BP"
$ws1.Range("C6").Value = "chondroitin sulfate proteoglycan metabolic process"
$ws1.Range("A7").Value = "GO:0080134"
$ws1.Range("B7").Value = "GO:BP"
$ws1.Range("C7").Value = "regulation of response to stress"
$ws1.Range("A8").Value = "GO:1902882"
$ws1.Range("B8").Value = "GO:BP"
$ws1.Range("C8").Value = "regulation of response to oxidative stress"
$ws1.Range("A9").Value = "GO:1902884"
$ws1.Range("B9").Value = "GO:BP"
$ws1.Range("C9").Value = "positive regulation of response to oxidative stress"
$ws1.Range("A10").Value = "GO:0006023"
$ws1.Range("B10").Value = "GO:BP"
$ws1.Range("C10").Value = "aminoglycan biosynthetic process"
$ws1.Range("A11").Value = "GO:0006024"
$ws1.Range("B11").Value = "GO:BP"
$ws1.Range("C11").Value = "glycosaminoglycan biosynthetic process"
$ws1.Range("A12").Value = "GO:0030203"
$ws1.Range("B12").Value = "GO:BP"
$ws1.Range("C12").Value = "glycosaminoglycan metabolic process"
$ws1.Range("A13").Value = "GO:0030166"
$ws1.Range("B13").Value = "GO:BP"
$ws1.Range("C13").Value = "proteoglycan biosynthetic process"
$ws1.Range("A14").Value = "GO:0006029"
$ws1.Range("B14").Value = "GO:BP"
$ws1.Range("C14").Value = "proteoglycan metabolic process"
$ws1.Range("A15").Value = "GO:0006022"
$ws1.Range("B15").Value = "GO:BP"
$ws1.Range("C15").Value = "aminoglycan metabolic process"
$ws1.Range("A16").Value = "GO:0048584"
$ws1.Range("B16").Value = "GO:BP"
$ws1.Range("C16").Value = "positive regulation of response to stimulus"
$ws1.Range("A17").Value = "GO:0007606"
$ws1.Range("B17").Value = "GO:BP"
$ws1.Range("C17").Value = "sensory perception of chemical stimulus"
$ws1.Range("A18").Value = "GO:0007600"
$ws1.Range("B18").Value = "GO:BP"
$ws1.Range("C18").Value = "sensory perception"
$ws1.Range("A19").Value = "GO:0006979"
$ws1.Range("B19").Value = "GO:BP"
$ws1.Range("C19").Value = "response to oxidative stress"
$ws1.Range("A20").Value = "GO:0009101"
$ws1.Range("B20").Value = "GO:BP"
$ws1.Range("C20").Value = "glycoprotein biosynthetic process"
$ws1.Range("A21").Value = "GO:0009100"
$ws1.Range("B21").Value = "GO:BP"
$ws1.Range("C21").Value = "glycoprotein metabolic process"
$ws1.Range("A22").Value = "GO:0044272"
$ws1.Range("B22").Value = "GO:BP"
$ws1.Range("C22").Value = "sulfur compound biosynthetic process"
$ws1.Range("A23").Value = "GO:0050877"
$ws1.Range("B23").Value = "GO:BP"
$ws1.Range("C23").Value = "nervous system process"
$ws1.Range("A24").Value = "GO:0003008"
$ws1.Range("B24").Value = "GO:BP"
$ws1.Range("C24").Value = "system process"
$ws1.Range("A25").Value = "GO:0006790"
$ws1.Range("B25").Value = "GO:BP"
$ws1.Range("C25").Value = "sulfur compound metabolic process"
$ws1.Range("A26").Value = "GO:1901137"
$ws1.Range("B26").Value = "GO:BP"
$ws1.Range("C26").Value = "carbohydrate derivative biosynthetic process"
$ws1.Range("A27").Value = "GO:1901135"
$ws1.Range("B27").Value = "GO:BP"
$ws1.Range("C27").Value = "carbohydrate derivative metabolic process"
$ws1.Range("A28").Value = "GO:0016021"
$ws1.Range("B28").Value = "GO:CC"
$ws1.Range("C28").Value = "integral component of membrane"
$ws1.Range("A29").Value = "GO:0031224"
$ws1.Range("B29").Value = "GO:CC"
$ws1.Range("C29").Value = "intrinsic component of membrane"
$ws1.Range("A30").Value = "GO:0016020"
$ws1.Range("B30").Value = "GO:CC"
$ws1.Range("C30").Value = "membrane"
$ws1.Range("A31").Value = "GO:0008146"
$ws1.Range("B31").Value = "GO:MF"
$ws1.Range("C31").Value = "sulfotransferase activity"
$ws1.Range("A32").Value = "GO:0034481"
$ws1.Range("B32").Value = "GO:MF"
$ws1.Range("C32").Value = "chondroitin sulfotransferase activity"
$ws1.Range("A33").Value = "GO:0047756"
$ws1.Range("B33").Value = "GO:MF"
$ws1.Range("C33").Value = "chondroitin 4-sulfotransferase activity"
$ws1.Range("A34").Value = "GO:0016782"
$ws1.Range("B34").Value = "GO:MF"
$ws1.Range("C34").Value = "transferase activity, transferring sulfur-containing groups"

# --- GC-rich genes ---
$ws2.Range("A2").Value = "GO:0043043"
$ws2.Range("B2").Value = "GO:BP"
$ws2.Range("C2").Value = "peptide biosynthetic process"
$ws2.Range("A3").Value = "GO:0006412"
$ws2.Range("B3").Value = "GO:BP"
$ws2.Range("C3").Value = "translation"
$ws2.Range("A4").Value = "GO:0043604"
$ws2.Range("B4").Value = "GO:BP"
$ws2.Range("C4").Value = "amide biosynthetic process"
$ws2.Range("A5").Value = "GO:0006518"
$ws2.Range("B5").Value = "GO:BP"
$ws2.Range("C5").Value = "peptide metabolic process"
$ws2.Range("A6").Value = "GO:0043603"
$ws2.Range("B6").Value = "GO:BP"
$ws2.Range("C6").Value = "cellular amide metabolic process"
$ws2.Range("A7").Value = "GO:1901566"
$ws2.Range("B7").Value = "GO:BP"
$ws2.Range("C7").Value = "organonitrogen compound biosynthetic process"
$ws2.Range("A8").Value = "GO:0044271"
$ws2.Range("B8").Value = "GO:BP"
$ws2.Range("C8").Value = "cellular nitrogen compound biosynthetic process"
$ws2.Range("A9").Value = "GO:0010467"
$ws2.Range("B9").Value = "GO:BP"
$ws2.Range("C9").Value = "gene expression"
$ws2.Range("A10").Value = "GO:0034641"
$ws2.Range("B10").Value = "GO:BP"
$ws2.Range("C10").Value = "cellular nitrogen compound metabolic process"
$ws2.Range("A11").Value = "GO:0034645"
$ws2.Range("B11").Value = "GO:BP"
$ws2.Range("C11").Value = "cellular macromolecule biosynthetic process"
$ws2.Range("A12").Value = "GO:0044249"
$ws2.Range("B12").Value = "GO:BP"
$ws2.Range("C12").Value = "cellular biosynthetic process"
$ws2.Range("A13").Value = "GO:0009059"
$ws2.Range("B13").Value = "GO:BP"
$ws2.Range("C13").Value = "macromolecule biosynthetic process"
$ws2.Range("A14").Value = "GO:1901576"
$ws2.Range("B14").Value = "GO:BP"
$ws2.Range("C14").Value = "organic substance biosynthetic process"
$ws2.Range("A15").Value = "GO:0009058"
$ws2.Range("B15").Value = "GO:BP"
$ws2.Range("C15").Value = "biosynthetic process"
$ws2.Range("A16").Value = "GO:0044267"
$ws2.Range("B16").Value = "GO:BP"
$ws2.Range("C16").Value = "cellular protein metabolic process"
$ws2.Range("A17").Value = "GO:0043232"
$ws2.Range("B17").Value = "GO:CC"
$ws2.Range("C17").Value = "intracellular non-membrane-bounded organelle"
$ws2.Range("A18").Value = "GO:0043228"
$ws2.Range("B18").Value = "GO:CC"
$ws2.Range("C18").Value = "non-membrane-bounded organelle"
$ws2.Range("A19").Value = "GO:0044815"
$ws2.Range("B19").Value = "GO:CC"
$ws2.Range("C19").Value = "DNA packaging complex"
$ws2.Range("A20").Value = "GO:0032993"
$ws2.Range("B20").Value = "GO:CC"
$ws2.Range("C20").Value = "protein-DNA complex"
$ws2.Range("A21").Value = "GO:0000786"
$ws2.Range("B21").Value = "GO:CC"
$ws2.Range("C21").Value = "nucleosome"
$ws2.Range("A22").Value = "GO:0005840"
$ws2.Range("B22").Value = "GO:CC"
$ws2.Range("C22").Value = "ribosome"
$ws2.Range("A23").Value = "GO:0000785"
$ws2.Range("B23").Value = "GO:CC"
$ws2.Range("C23").Value = "chromatin"
$ws2.Range("A24").Value = "GO:0005694"
$ws2.Range("B24").Value = "GO:CC"
$ws2.Range("C24").Value = "chromosome"
$ws2.Range("A25").Value = "GO:0042302"
$ws2.Range("B25").Value = "GO:MF"
$ws2.Range("C25").Value = "structural constituent of cuticle"
$ws2.Range("A26").Value = "GO:0005198"
$ws2.Range("B26").Value = "GO:MF"
$ws2.Range("C26").Value = "structural molecule activity"
$ws2.Range("A27").Value = "GO:0003735"
$ws2.Range("B27").Value = "GO:MF"
$ws2.Range("C27").Value = "structural constituent of ribosome"
$ws2.Range("A28").Value = "GO:0046982"
$ws2.Range("B28").Value = "GO:MF"
$ws2.Range("C28").Value = "protein heterodimerization activity"
$ws2.Range("A29").Value = "GO:0046983"
$ws2.Range("B29").Value = "GO:MF"
$ws2.Range("C29").Value = "protein dimerization activity"

# --- Least-adapted genes ---
$ws3.Range("A2").Value = "GO:0015074"
$ws3.Range("B2").Value = "GO:BP"
$ws3.Range("C2").Value = "DNA integration"
$ws3.Range("A3").Value = "GO:0003676"
$ws3.Range("B3").Value = "GO:MF"
$ws3.Range("C3").Value = "nucleic acid binding"
$ws3.Range("A4").Value = "GO:0006259"
$ws3.Range("B4").Value = "GO:BP"
$ws3.Range("C4").Value = "DNA metabolic process"
$ws3.Range("A5").Value = "GO:0090304"
$ws3.Range("B5").Value = "GO:BP"
$ws3.Range("C5").Value = "nucleic acid metabolic process"

# --- Most-adapted genes ---
$ws4.Range("A2").Value = "GO:0005198"
$ws4.Range("B2").Value = "GO:MF"
$ws4.Range("C2").Value = "structural molecule activity"
$ws4.Range("A3").Value = "GO:0043232"
$ws4.Range("B3").Value = "GO:CC"
$ws4.Range("C3").Value = "intracellular non-membrane-bounded organelle"
$ws4.Range("A4").Value = "GO:0043228"
$ws4.Range("B4").Value = "GO:CC"
$ws4.Range("C4").Value = "non-membrane-bounded organelle"
$ws4.Range("A5").Value = "GO:0042302"
$ws4.Range("B5").Value = "GO:MF"
$ws4.Range("C5").Value = "structural constituent of cuticle"
$ws4.Range("A6").Value = "GO:0005840"
$ws4.Range("B6").Value = "GO:CC"
$ws4.Range("C6").Value = "ribosome"
$ws4.Range("A7").Value = "GO:0046034"
$ws4.Range("B7").Value = "GO:BP"
$ws4.Range("C7").Value = "ATP metabolic process"
